$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1368.2659
$ws.Range("J17").Value = 1368.2659
$ws.Range("L17").Value = 4104.7977
$ws.Range("N17").Value = -4440.7977
$ws.Range("H88").Value = 7320
$ws.Range("I88").Value = 7557.143
$ws.Range("J88").Value = 6490
$ws.Range("K88").Value = 7557.143
$ws.Range("L88").Value = 6490
$ws.Range("M88").Value = -7151.143
$ws.Range("N88").Value = -7302
$ws.Range("H91").Value = 7320
$ws.Range("I91").Value = 7557.143
$ws.Range("J91").Value = 6490
$ws.Range("K91").Value = 7557.143
$ws.Range("L91").Value = 6490
$ws.Range("M91").Value = -6153.143
$ws.Range("N91").Value = -9298
$ws.Range("H112").Value = 2139.3333
$ws.Range("I112").Value = 300
$ws.Range("J112").Value = 2332.9473
$ws.Range("K112").Value = 900
$ws.Range("L112").Value = 6998.841899999999
$ws.Range("M112").Value = 208
$ws.Range("N112").Value = -9214.841899999999
$ws.Range("H116").Value = 6298272
$ws.Range("I116").Value = 6747991.5
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 6747991.5
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = -6744549.5
$ws.Range("N116").Value = -9084
$ws.Range("H128").Value = 9556
$ws.Range("J128").Value = 9556
$ws.Range("L128").Value = 9556
$ws.Range("N128").Value = -19516
$ws.Range("H138").Value = 140826.08
$ws.Range("I138").Value = 455308.1
$ws.Range("J138").Value = 2454
$ws.Range("K138").Value = 1365924.3
$ws.Range("L138").Value = 7362
$ws.Range("M138").Value = -1360784.3
$ws.Range("N138").Value = -17642

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5321278.5
$ws.Range("I74").Value = 10002335
$ws.Range("J74").Value = 1895.409
$ws.Range("K74").Value = 10002335
$ws.Range("L74").Value = 1895.409
$ws.Range("M74").Value = -10001461
$ws.Range("N74").Value = -3643.409
$ws.Range("H77").Value = 5321278.5
$ws.Range("I77").Value = 10002335
$ws.Range("J77").Value = 1895.409
$ws.Range("K77").Value = 50011675
$ws.Range("L77").Value = 9477.045
$ws.Range("M77").Value = -50007307
$ws.Range("N77").Value = -18213.045
$ws.Range("H102").Value = 50001496
$ws.Range("I102").Value = 111112110
$ws.Range("J102").Value = 1906.1818
$ws.Range("K102").Value = 111112110
$ws.Range("L102").Value = 1906.1818
$ws.Range("M102").Value = -111110488
$ws.Range("N102").Value = -5150.1818
$ws.Range("H105").Value = 38000
$ws.Range("J105").Value = 38000
$ws.Range("L105").Value = 38000
$ws.Range("N105").Value = -44988
$ws.Range("H110").Value = 8807.125
$ws.Range("I110").Value = 12442.1
$ws.Range("J110").Value = 2748.8333
$ws.Range("K110").Value = 12442.1
$ws.Range("L110").Value = 2748.8333
$ws.Range("M110").Value = -10397.1
$ws.Range("N110").Value = -6838.8333
$ws.Range("H132").Value = 2290.1538
$ws.Range("I132").Value = 2049.8096
$ws.Range("J132").Value = 3299.6
$ws.Range("K132").Value = 6149.4288
$ws.Range("L132").Value = 9898.799999999999
$ws.Range("M132").Value = -3619.4288
$ws.Range("N132").Value = -14958.8

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 80684
$ws.Range("J42").Value = 80684
$ws.Range("L42").Value = 80684
$ws.Range("N42").Value = -81340
$ws.Range("H47").Value = 80684
$ws.Range("J47").Value = 80684
$ws.Range("L47").Value = 80684
$ws.Range("N47").Value = -81724
$ws.Range("H104").Value = 29800
$ws.Range("J104").Value = 29800
$ws.Range("L104").Value = 29800
$ws.Range("N104").Value = -36788
$ws.Range("H105").Value = 2424.85
$ws.Range("I105").Value = 2198.182
$ws.Range("J105").Value = 2701.889
$ws.Range("K105").Value = 2198.182
$ws.Range("L105").Value = 2701.889
$ws.Range("M105").Value = -451.1819999999998
$ws.Range("N105").Value = -6195.889

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9549080
$ws.Range("I31").Value = 7731794.5
$ws.Range("J31").Value = 12502169
$ws.Range("K31").Value = 7731794.5
$ws.Range("L31").Value = 12502169
$ws.Range("M31").Value = -7731499.5
$ws.Range("N31").Value = -12502759
$ws.Range("H34").Value = 9549080
$ws.Range("I34").Value = 7731794.5
$ws.Range("J34").Value = 12502169
$ws.Range("K34").Value = 7731794.5
$ws.Range("L34").Value = 12502169
$ws.Range("M34").Value = -7731592.5
$ws.Range("N34").Value = -12502573
$ws.Range("H107").Value = 1188.8077
$ws.Range("I107").Value = 1154.2941
$ws.Range("J107").Value = 1254
$ws.Range("K107").Value = 1154.2941
$ws.Range("L107").Value = 1254
$ws.Range("M107").Value = 765.7058999999999
$ws.Range("N107").Value = -5094

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1553.6875
$ws.Range("I80").Value = 900
$ws.Range("J80").Value = 1597.2667
$ws.Range("K80").Value = 2700
$ws.Range("L80").Value = 4791.800099999999
$ws.Range("M80").Value = -1764
$ws.Range("N80").Value = -6663.800099999999
$ws.Range("H83").Value = 1553.6875
$ws.Range("I83").Value = 900
$ws.Range("J83").Value = 1597.2667
$ws.Range("K83").Value = 8100
$ws.Range("L83").Value = 14375.4003
$ws.Range("M83").Value = -3420
$ws.Range("N83").Value = -23735.4003
$ws.Range("H104").Value = 4092
$ws.Range("J104").Value = 4040.6667
$ws.Range("L104").Value = 12122.0001
$ws.Range("N104").Value = -17364.0001
$ws.Range("H131").Value = 3367868.5
$ws.Range("I131").Value = 419.47058
$ws.Range("J131").Value = 4065998.2
$ws.Range("K131").Value = 1258.41174
$ws.Range("L131").Value = 12197994.6
$ws.Range("M131").Value = 3781.58826
$ws.Range("N131").Value = -12208074.6

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2446.3635
$ws.Range("I102").Value = 2434.4443
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 2434.4443
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -812.4443000000001
$ws.Range("N102").Value = -5744
$ws.Range("H113").Value = 2082.4
$ws.Range("I113").Value = 1818.5
$ws.Range("K113").Value = 1818.5
$ws.Range("M113").Value = 351.5
$ws.Range("H122").Value = 3671.818
$ws.Range("I122").Value = 2706.818
$ws.Range("J122").Value = 4636.8184
$ws.Range("K122").Value = 8120.454000000001
$ws.Range("L122").Value = 13910.4552
$ws.Range("M122").Value = -5670.454000000001
$ws.Range("N122").Value = -18810.4552

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 11364942
$ws.Range("I93").Value = 13159178
$ws.Range("J93").Value = 1444.1666
$ws.Range("K93").Value = 13159178
$ws.Range("L93").Value = 1444.1666
$ws.Range("M93").Value = -13157930
$ws.Range("N93").Value = -3940.1666
$ws.Range("H132").Value = 11256557
$ws.Range("I132").Value = 22510314
$ws.Range("J132").Value = 2799.72
$ws.Range("K132").Value = 67530942
$ws.Range("L132").Value = 8399.16
$ws.Range("M132").Value = -67528412
$ws.Range("N132").Value = -13459.16

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 35000
$ws.Range("L106").Value = 35000
